$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 9: new date and new volume
$ws.Range("D9").Value = 44769
$ws.Range("J9").Value = 50

# Append new row 10 (copy of the former row 9 data, pre-edit)
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Vega Modelo de Temuco"
$ws.Range("C10").Value = "La Araucanía"
$ws.Range("D10").Value = 44757
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 100112036
$ws.Range("G10").Value = "Caigua"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 30
$ws.Range("K10").Value = 20000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 20000
$ws.Range("N10").Value = "$/caja 15 kilos"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 1333
$ws.Range("Q10").Value = 15
$ws.Range("R10").Value = "Hortaliza"
